# chore: update Sheets via scheduled runner
# Applies the price/profit refresh from the scraped diff across the FFXIV
# leve-profit workbook's per-job sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR).
# Each row's currentAveragePrice* (H/I/J/K/L) and derived profit (M/N) columns
# are updated to the new market-board snapshot values; a handful of cells
# are newly populated or cleared entirely to mirror sparse-cell add/remove
# in the source diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 47.6  # was 56.4
$ws.Range("I11").Value = 47.6  # was 56.4
$ws.Range("K11").Value = 47.6  # was 56.4
$ws.Range("M11").Value = 92.40000000000001  # was 83.59999999999999
$ws.Range("H12").Value = 0  # was 234
$ws.Range("J12").Value = 0  # was 234
$ws.Range("L12").Value = 0  # was 234
$ws.Range("N12").ClearContents()  # was -574
$ws.Range("H18").Value = 2062.5  # was 1785.7142
$ws.Range("J18").Value = 4000  # was 0
$ws.Range("L18").Value = 4000  # was 0
$ws.Range("N18").Value = -4568  # was None
$ws.Range("H38").Value = 1168  # was 626
$ws.Range("I38").Value = 1168  # was 501.33334
$ws.Range("J38").Value = 0  # was 1000
$ws.Range("K38").Value = 3504  # was 1504.00002
$ws.Range("L38").Value = 0  # was 3000
$ws.Range("M38").Value = -3132  # was -1132.00002
$ws.Range("N38").ClearContents()  # was -3744
$ws.Range("H40").Value = 2333.3333  # was 2000
$ws.Range("J40").Value = 1500  # was 1333.3334
$ws.Range("L40").Value = 1500  # was 1333.3334
$ws.Range("N40").Value = -1850  # was -1683.3334
$ws.Range("H58").Value = 5999.75  # was 7500
$ws.Range("J58").Value = 5999.75  # was 7500
$ws.Range("L58").Value = 17999.25  # was 22500
$ws.Range("N58").Value = -18299.25  # was -22800
$ws.Range("H87").Value = 99354  # was 98853.5
$ws.Range("J87").Value = 99354  # was 98853.5
$ws.Range("L87").Value = 99354  # was 98853.5
$ws.Range("N87").Value = -101850  # was -101349.5
$ws.Range("H90").Value = 99354  # was 98853.5
$ws.Range("J90").Value = 99354  # was 98853.5
$ws.Range("L90").Value = 298062  # was 296560.5
$ws.Range("N90").Value = -310542  # was -309040.5
$ws.Range("H100").Value = 4620  # was 6666.6665
$ws.Range("I100").Value = 4326.6665  # was 5000
$ws.Range("J100").Value = 5500  # was 10000
$ws.Range("K100").Value = 4326.6665  # was 5000
$ws.Range("L100").Value = 5500  # was 10000
$ws.Range("M100").Value = -3785.6665  # was -4459
$ws.Range("N100").Value = -6582  # was -11082
$ws.Range("H111").Value = 254.75  # was 3942.25
$ws.Range("I111").Value = 210  # was 5156.3335
$ws.Range("J111").Value = 299.5  # was 300
$ws.Range("K111").Value = 630  # was 15469.0005
$ws.Range("L111").Value = 898.5  # was 900
$ws.Range("M111").Value = 2437  # was -12402.0005
$ws.Range("N111").Value = -7032.5  # was -7034
$ws.Range("H135").Value = 514.36365  # was 462.2
$ws.Range("I135").Value = 514.36365  # was 462.2
$ws.Range("K135").Value = 4629.27285  # was 4159.8
$ws.Range("M135").Value = -2094.27285  # was -1624.8
$ws.Range("H138").Value = 10665.866  # was 11641.857
$ws.Range("I138").Value = 999.4  # was 999.75
$ws.Range("J138").Value = 15499.1  # was 15898.7
$ws.Range("K138").Value = 2998.2  # was 2999.25
$ws.Range("L138").Value = 46497.3  # was 47696.10000000001
$ws.Range("M138").Value = 2141.8  # was 2140.75
$ws.Range("N138").Value = -56777.3  # was -57976.10000000001
$ws.Range("H141").Value = 4999.75  # was 5666.3335
$ws.Range("I141").Value = 4999.75  # was 5666.3335
$ws.Range("K141").Value = 14999.25  # was 16999.0005
$ws.Range("M141").Value = -9819.25  # was -11819.0005

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4450.8  # was 4810
$ws.Range("I61").Value = 4450.8  # was 4810
$ws.Range("K61").Value = 4450.8  # was 4810
$ws.Range("M61").Value = -4238.8  # was -4598
$ws.Range("H102").Value = 1673.0714  # was 1983.909
$ws.Range("I102").Value = 1673.0714  # was 1983.909
$ws.Range("K102").Value = 1673.0714  # was 1983.909
$ws.Range("M102").Value = -51.07140000000004  # was -361.9090000000001
$ws.Range("H113").Value = 0  # was 20000
$ws.Range("J113").Value = 0  # was 20000
$ws.Range("L113").Value = 0  # was 20000
$ws.Range("N113").ClearContents()  # was -28678
$ws.Range("H136").Value = 4450.8  # was 4810
$ws.Range("I136").Value = 4450.8  # was 4810
$ws.Range("K136").Value = 13352.4  # was 14430
$ws.Range("M136").Value = -10802.4  # was -11880

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 10559.7  # was 7995.7144
$ws.Range("I86").Value = 3819.8  # was 3185.1428
$ws.Range("J86").Value = 17299.6  # was 12806.286
$ws.Range("K86").Value = 3819.8  # was 3185.1428
$ws.Range("L86").Value = 17299.6  # was 12806.286
$ws.Range("M86").Value = -2696.8  # was -2062.1428
$ws.Range("N86").Value = -19545.6  # was -15052.286
$ws.Range("H89").Value = 10559.7  # was 7995.7144
$ws.Range("I89").Value = 3819.8  # was 3185.1428
$ws.Range("J89").Value = 17299.6  # was 12806.286
$ws.Range("K89").Value = 19099  # was 15925.714
$ws.Range("L89").Value = 86498  # was 64031.43
$ws.Range("M89").Value = -13483  # was -10309.714
$ws.Range("N89").Value = -97730  # was -75263.42999999999
$ws.Range("H105").Value = 11200.857  # was 11200.714
$ws.Range("I105").Value = 12484.333  # was 11200.714
$ws.Range("J105").Value = 3500  # was 0
$ws.Range("K105").Value = 12484.333  # was 11200.714
$ws.Range("L105").Value = 3500  # was 0
$ws.Range("M105").Value = -10737.333  # was -9453.714
$ws.Range("N105").Value = -6994  # was None
$ws.Range("H134").Value = 4000  # was 1999.6666
$ws.Range("I134").Value = 4000  # was 1999.6666
$ws.Range("K134").Value = 12000  # was 5998.9998
$ws.Range("M134").Value = -9465  # was -3463.9998

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 33982  # was 17999
$ws.Range("J41").Value = 49965  # was 0
$ws.Range("L41").Value = 49965  # was 0
$ws.Range("N41").Value = -50821  # was None
$ws.Range("H51").Value = 52030.4  # was 52540.75
$ws.Range("J51").Value = 52030.4  # was 52540.75
$ws.Range("L51").Value = 52030.4  # was 52540.75
$ws.Range("N51").Value = -53502.4  # was -54012.75
$ws.Range("H60").Value = 39697.668  # was 34995
$ws.Range("J60").Value = 49549  # was 49995
$ws.Range("L60").Value = 49549  # was 49995
$ws.Range("N60").Value = -50571  # was -51017
$ws.Range("H61").Value = 52030.4  # was 52540.75
$ws.Range("J61").Value = 52030.4  # was 52540.75
$ws.Range("L61").Value = 52030.4  # was 52540.75
$ws.Range("N61").Value = -52726.4  # was -53236.75
$ws.Range("H86").Value = 19498.166  # was 17138.285
$ws.Range("I86").Value = 27748.5  # was 22794.6
$ws.Range("K86").Value = 27748.5  # was 22794.6
$ws.Range("M86").Value = -26625.5  # was -21671.6
$ws.Range("H89").Value = 19498.166  # was 17138.285
$ws.Range("I89").Value = 27748.5  # was 22794.6
$ws.Range("K89").Value = 138742.5  # was 113973
$ws.Range("M89").Value = -133126.5  # was -108357
$ws.Range("H105").Value = 3000  # was 2418
$ws.Range("I105").Value = 3000  # was 2611.6
$ws.Range("J105").Value = 0  # was 1450
$ws.Range("K105").Value = 3000  # was 2611.6
$ws.Range("L105").Value = 0  # was 1450
$ws.Range("M105").Value = -1253  # was -864.5999999999999
$ws.Range("N105").ClearContents()  # was -4944

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 100344  # was 0
$ws.Range("J95").Value = 100344  # was 0
$ws.Range("L95").Value = 100344  # was 0
$ws.Range("N95").Value = -105836  # was None
$ws.Range("H102").Value = 774.5  # was 777
$ws.Range("I102").Value = 774.5  # was 777
$ws.Range("K102").Value = 774.5  # was 777
$ws.Range("M102").Value = 847.5  # was 845
$ws.Range("H107").Value = 819.8  # was 820
$ws.Range("I107").Value = 1599.5  # was 1600
$ws.Range("K107").Value = 1599.5  # was 1600
$ws.Range("M107").Value = 320.5  # was 320
$ws.Range("H113").Value = 2126.4167  # was 2183.5833
$ws.Range("I113").Value = 2361.8  # was 2301.3333
$ws.Range("J113").Value = 1958.2858  # was 2065.8333
$ws.Range("K113").Value = 2361.8  # was 2301.3333
$ws.Range("L113").Value = 1958.2858  # was 2065.8333
$ws.Range("M113").Value = -191.8000000000002  # was -131.3332999999998
$ws.Range("N113").Value = -6298.2858  # was -6405.8333
$ws.Range("H132").Value = 4723.3335  # was 5249.6665
$ws.Range("I132").Value = 3802.4  # was 4000
$ws.Range("K132").Value = 11407.2  # was 12000
$ws.Range("M132").Value = -8877.200000000001  # was -9470

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7639.722  # was 6960.75
$ws.Range("J22").Value = 10000  # was 7385.7144
$ws.Range("L22").Value = 10000  # was 7385.7144
$ws.Range("N22").Value = -10590  # was -7975.7144
$ws.Range("H27").Value = 7639.722  # was 6960.75
$ws.Range("J27").Value = 10000  # was 7385.7144
$ws.Range("L27").Value = 10000  # was 7385.7144
$ws.Range("N27").Value = -10214  # was -7599.7144
$ws.Range("H46").Value = 3216.5  # was 3360
$ws.Range("I46").Value = 3074.75  # was 3266.6667
$ws.Range("K46").Value = 3074.75  # was 3266.6667
$ws.Range("M46").Value = -2886.75  # was -3078.6667
$ws.Range("H68").Value = 1773  # was 1676.7142
$ws.Range("I68").Value = 1773  # was 1676.7142
$ws.Range("K68").Value = 1773  # was 1676.7142
$ws.Range("M68").Value = -1024  # was -927.7141999999999
$ws.Range("H71").Value = 1773  # was 1676.7142
$ws.Range("I71").Value = 1773  # was 1676.7142
$ws.Range("K71").Value = 8865  # was 8383.571
$ws.Range("M71").Value = -5121  # was -4639.571
$ws.Range("H93").Value = 0  # was 1902
$ws.Range("I93").Value = 0  # was 1902
$ws.Range("K93").Value = 0  # was 1902
$ws.Range("M93").ClearContents()  # was -654

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2897.5  # was 1197.5
$ws.Range("I122").Value = 2897.5  # was 1197.5
$ws.Range("K122").Value = 8692.5  # was 3592.5
$ws.Range("M122").Value = -6242.5  # was -1142.5
$ws.Range("H126").Value = 4566.6665  # was 3628.5715
$ws.Range("I126").Value = 3880  # was 3628.5715
$ws.Range("J126").Value = 8000  # was 0
$ws.Range("K126").Value = 11640  # was 10885.7145
$ws.Range("L126").Value = 24000  # was 0
$ws.Range("M126").Value = -9170  # was -8415.7145
$ws.Range("N126").Value = -28940  # was None

